$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new control-reference values for rows 3 and 2 (in this order so the
# shared-strings table ends up with the same ordering as the saved workbook)
$ws.Range("E3").Value = "AC-11, AC-12"
$ws.Range("E2").Value = "IA-2"

# Update the header of column E: "Controls Reference" -> "Controls_Reference"
$ws.Range("E1").Value = "Controls_Reference"

# Move the active selection to E1 (matches the saved selection in the diff)
$ws.Range("E1").Select()
